$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new PASS! status and timestamp to row 2 (Status / Time columns)
$ws.Range("G2").Value = "PASS!"
$ws.Range("H2").Value = "03:16 PM"

# Update the active cell selection to H2
$ws.Range("H2").Select()
